$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the
#    document's title (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$labelStart = $metaRange.Start
$label = "Meta description"
$metaRange.InsertAfter($label + ": Read our review of Dragon Spin Pick n Mix, an online slot game with 5 bonuses, a Dragon Pot Bonus, and a high RTP. Play it for free or real money.")

$labelRange = $d.Range($labelStart, $labelStart + $label.Length)
$labelRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicate bold title paragraph that used to sit right
#    before the italic "Read our review..." paragraph near the end.
# ------------------------------------------------------------------
$oldTitle = "Play Dragon Spin Pick n Mix Free Slot Game | Review 2021"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd()
    if ($i -gt 1 -and $text -eq $oldTitle) {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) Replace the text of the remaining italic paragraph (formerly the
#    meta description) with the new image-prompt text, keeping its
#    italic formatting intact.
# ------------------------------------------------------------------
$oldBody = "Read our review of Dragon Spin Pick n Mix, an online slot game with 5 bonuses, a Dragon Pot Bonus, and a high RTP. Play it for free or real money."
$newBody = "Create a feature image for Dragon Spin Pick n Mix that features a happy Maya warrior with glasses in a cartoon-style. The warrior should be surrounded by dragons and treasure, creating a sense of adventure and excitement. The image should be bright and colorful, with an Eastern-inspired design to match the game's theme. The Maya warrior should be depicted as confident and triumphant, holding a winning jackpot symbol above their head. This will appeal to the slot game's adventurous and adventurous players, while also highlighting the exciting gameplay and potential for big wins."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd()
    if ($text -eq $oldBody) {
        $r = $para.Range
        $bodyRange = $d.Range($r.Start, $r.End - 1)
        $bodyRange.Text = $newBody
        break
    }
}
